$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 181.4944075
$ws.Range("H2").Value = 362.988815
$ws.Range("I2").Value = 0.2239486468210351
$ws.Range("J2").Value = 0.1654349085470023
$ws.Range("O2").Value = 0.1656403066315349
$ws.Range("P2").Value = 0.2294568116288535
$ws.Range("Q2").Value = 19.50841037522417
$ws.Range("R2").Value = 117.050462251345
$ws.Range("S2").Value = 0.03709492252915356
$ws.Range("T2").Value = 0.03796016664730612
$ws.Range("G3").Value = 181.4944075
$ws.Range("H3").Value = 362.988815
$ws.Range("I3").Value = 0.2239486468210351
$ws.Range("J3").Value = 0.1654349085470023
$ws.Range("M3").Value = 0.5414345
$ws.Range("N3").Value = 1.082869
$ws.Range("O3").Value = 0.8343596933684652
$ws.Range("P3").Value = 0.7705431883711465
$ws.Range("Q3").Value = 98.26733377755878
$ws.Range("R3").Value = 393.0693351102351
$ws.Range("S3").Value = 0.1868537242918815
$ws.Range("T3").Value = 0.1274747418996962
$ws.Range("I4").Value = 0.07700606288633029
$ws.Range("J4").Value = 0.08532865336765341
$ws.Range("O4").Value = 0.1656403066315349
$ws.Range("P4").Value = 0.2294568116288535
$ws.Range("S4").Value = 0.01275530786897901
$ws.Range("T4").Value = 0.01957924074232539
$ws.Range("I5").Value = 0.07700606288633029
$ws.Range("J5").Value = 0.08532865336765341
$ws.Range("M5").Value = 0.5414345
$ws.Range("N5").Value = 1.082869
$ws.Range("O5").Value = 0.8343596933684652
$ws.Range("P5").Value = 0.7705431883711465
$ws.Range("Q5").Value = 33.78980222458717
$ws.Range("R5").Value = 202.738813347523
$ws.Range("S5").Value = 0.06425075501735128
$ws.Range("T5").Value = 0.06574941262532803
$ws.Range("G6").Value = 171.9980316666667
$ws.Range("H6").Value = 515.994095
$ws.Range("I6").Value = 0.2122309275432167
$ws.Range("J6").Value = 0.235168226649403
$ws.Range("O6").Value = 0.1656403066315349
$ws.Range("P6").Value = 0.2294568116288535
$ws.Range("Q6").Value = 18.48766709510944
$ws.Range("R6").Value = 166.389003855985
$ws.Range("S6").Value = 0.03515399591495347
$ws.Range("T6").Value = 0.0539609514833836
$ws.Range("G7").Value = 171.9980316666667
$ws.Range("H7").Value = 515.994095
$ws.Range("I7").Value = 0.2122309275432167
$ws.Range("J7").Value = 0.235168226649403
$ws.Range("M7").Value = 0.5414345
$ws.Range("N7").Value = 1.082869
$ws.Range("O7").Value = 0.8343596933684652
$ws.Range("P7").Value = 0.7705431883711465
$ws.Range("Q7").Value = 93.12566827642584
$ws.Range("R7").Value = 558.7540096585551
$ws.Range("S7").Value = 0.1770769316282632
$ws.Range("T7").Value = 0.1812072751660194
$ws.Range("G8").Value = 55.64279550000001
$ws.Range("H8").Value = 111.285591
$ws.Range("I8").Value = 0.06865847234198982
$ws.Range("J8").Value = 0.05071925307032974
$ws.Range("O8").Value = 0.1656403066315349
$ws.Range("P8").Value = 0.2294568116288535
$ws.Range("Q8").Value = 5.980914255105501
$ws.Range("R8").Value = 35.88548553063301
$ws.Range("S8").Value = 0.01137261041157995
$ws.Range("T8").Value = 0.0116378780977148
$ws.Range("G9").Value = 55.64279550000001
$ws.Range("H9").Value = 111.285591
$ws.Range("I9").Value = 0.06865847234198982
$ws.Range("J9").Value = 0.05071925307032974
$ws.Range("M9").Value = 0.5414345
$ws.Range("N9").Value = 1.082869
$ws.Range("O9").Value = 0.8343596933684652
$ws.Range("P9").Value = 0.7705431883711465
$ws.Range("Q9").Value = 30.12692916014475
$ws.Range("R9").Value = 120.507716640579
$ws.Range("S9").Value = 0.05728586193040987
$ws.Range("T9").Value = 0.03908137497261494
$ws.Range("G10").Value = 203.386317
$ws.Range("H10").Value = 610.158951
$ws.Range("I10").Value = 0.250961399315095
$ws.Range("J10").Value = 0.2780845747487284
$ws.Range("O10").Value = 0.1656403066315349
$ws.Range("P10").Value = 0.2294568116288535
$ws.Range("Q10").Value = 21.861520646257
$ws.Range("R10").Value = 196.753685816313
$ws.Range("S10").Value = 0.04156932313523141
$ws.Range("T10").Value = 0.06380839988500882
$ws.Range("G11").Value = 203.386317
$ws.Range("H11").Value = 610.158951
$ws.Range("I11").Value = 0.250961399315095
$ws.Range("J11").Value = 0.2780845747487284
$ws.Range("M11").Value = 0.5414345
$ws.Range("N11").Value = 1.082869
$ws.Range("O11").Value = 0.8343596933684652
$ws.Range("P11").Value = 0.7705431883711465
$ws.Range("Q11").Value = 110.1203688517365
$ws.Range("R11").Value = 660.722213110419
$ws.Range("S11").Value = 0.2093920761798636
$ws.Range("T11").Value = 0.2142761748637196
$ws.Range("G12").Value = 135.4992116666667
$ws.Range("H12").Value = 406.497635
$ws.Range("I12").Value = 0.167194491092333
$ws.Range("J12").Value = 0.1852643836168829
$ws.Range("O12").Value = 0.1656403066315349
$ws.Range("P12").Value = 0.2294568116288535
$ws.Range("Q12").Value = 14.56449409722278
$ws.Range("R12").Value = 131.080446875005
$ws.Range("S12").Value = 0.02769414677163747
$ws.Range("T12").Value = 0.04251017477311475
$ws.Range("G13").Value = 135.4992116666667
$ws.Range("H13").Value = 406.497635
$ws.Range("I13").Value = 0.167194491092333
$ws.Range("J13").Value = 0.1852643836168829
$ws.Range("M13").Value = 0.5414345
$ws.Range("N13").Value = 1.082869
$ws.Range("O13").Value = 0.8343596933684652
$ws.Range("P13").Value = 0.7705431883711465
$ws.Range("Q13").Value = 73.36394791913584
$ws.Range("R13").Value = 440.1836875148151
$ws.Range("S13").Value = 0.1395003443206956
$ws.Range("T13").Value = 0.1427542088437681
